# Fix: standardized_IS updated (Published)
# Split "Net Income (loss)" (row 22/23, shared the same label) into two
# distinct, correctly-labeled rows:
#   Row 22 -> Net Income (loss) (operations)                      [ProfitLoss]
#   Row 23 -> Net Income (loss) to parent (incl. Non contr. Interest) [NetIncomeLoss]
# (the underlying XBRL tag/description pair for the two rows is swapped
#  relative to the old sheet, since the consolidated "ProfitLoss" figure is
#  the operations-level number and "NetIncomeLoss" is the to-parent number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: "Net Income (loss)" -> "Net Income (loss) (operations)", tag -> ProfitLoss
$ws.Range("A22").Value = "Net Income (loss) (operations)"
$ws.Range("B22").Value = "ProfitLoss"
$ws.Range("C22").Value = "The consolidated profit or loss for the period, net of income taxes, including the portion attributable to the noncontrolling interest."

# --- Row 23: "Net Income (loss)" -> "Net Income (loss) to parent (incl. Non contr. Interest)", tag -> NetIncomeLoss
$ws.Range("A23").Value = "Net Income (loss) to parent (incl. Non contr. Interest)"
$ws.Range("B23").Value = "NetIncomeLoss"
$ws.Range("C23").Value = "The portion of profit or loss for the period, net of income taxes, which is attributable to the parent."

# Selection/active cell moved to A23 (matches author's final click position)
$ws.Range("A23").Select()
